$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Cell="E2"; Value=18.8715},
    @{Cell="F2"; Value=14.9419361744008},
    @{Cell="E3"; Value=36.57},
    @{Cell="F3"; Value=29.268438307697},
    @{Cell="I3"; Value=76.2393678027716},
    @{Cell="J3"; Value=2.89041931095767},
    @{Cell="E4"; Value=13.179},
    @{Cell="F4"; Value=9.86823214645805},
    @{Cell="I4"; Value=53.505899729772},
    @{Cell="J4"; Value=1.82286078684449},
    @{Cell="I5"; Value=65.1196514881153},
    @{Cell="J5"; Value=0.0430569909647724},
    @{Cell="I6"; Value=130.239146021488},
    @{Cell="J6"; Value=0.0114363542742864},
    @{Cell="E7"; Value=2.7485},
    @{Cell="F7"; Value=2.04554765899223},
    @{Cell="E8"; Value=1.1615},
    @{Cell="F8"; Value=0.998054133918157},
    @{Cell="I8"; Value=174.394640735725},
    @{Cell="J8"; Value=0.0513751201355547},
    @{Cell="E9"; Value=11.6955},
    @{Cell="F9"; Value=9.79555197301289},
    @{Cell="I9"; Value=95.7190811699057},
    @{Cell="J9"; Value=0.71593651049622},
    @{Cell="E10"; Value=10.3155},
    @{Cell="F10"; Value=8.02867368322654},
    @{Cell="E11"; Value=13.9265},
    @{Cell="F11"; Value=11.4118255281751},
    @{Cell="E12"; Value=2.116},
    @{Cell="F12"; Value=1.70138257245778},
    @{Cell="E13"; Value=6.5665},
    @{Cell="F13"; Value=5.07829538984966},
    @{Cell="I13"; Value=74.5496973329696},
    @{Cell="J13"; Value=0.453544554794711},
    @{Cell="E14"; Value=6.8885},
    @{Cell="F14"; Value=5.89147531695343},
    @{Cell="E15"; Value=4.7495},
    @{Cell="F15"; Value=4.00049963438626},
    @{Cell="E16"; Value=9.8785},
    @{Cell="F16"; Value=8.17480828450343},
    @{Cell="I16"; Value=87.420025653145},
    @{Cell="J16"; Value=0.615357115492503},
    @{Cell="E17"; Value=7.567},
    @{Cell="F17"; Value=6.21923221479076},
    @{Cell="I17"; Value=97.6440364315474},
    @{Cell="J17"; Value=0.486455245476295},
    @{Cell="E18"; Value=2.1045},
    @{Cell="F18"; Value=1.76770431707038},
    @{Cell="E19"; Value=0.3565},
    @{Cell="F19"; Value=0.306823999693611},
    @{Cell="I19"; Value=176.551365081164},
    @{Cell="J19"; Value=0.0661939008419496},
    @{Cell="E20"; Value=1.219},
    @{Cell="F20"; Value=0.906861891254617},
    @{Cell="I20"; Value=66.607300007954},
    @{Cell="J20"; Value=0.116616579068735},
    @{Cell="E22"; Value=8.7055},
    @{Cell="F22"; Value=5.94147543420233},
    @{Cell="I22"; Value=46.5946311448523},
    @{Cell="J22"; Value=1.67457427685377},
    @{Cell="E23"; Value=4.278},
    @{Cell="F23"; Value=3.45337595762792},
    @{Cell="E24"; Value=11.6725},
    @{Cell="F24"; Value=8.47692314013229},
    @{Cell="I24"; Value=52.9553438775132},
    @{Cell="J24"; Value=1.53451244378401},
    @{Cell="E25"; Value=17.664},
    @{Cell="F25"; Value=14.5839590397331},
    @{Cell="I25"; Value=92.7205132735387},
    @{Cell="J25"; Value=0.928616049320151},
    @{Cell="E26"; Value=22.011},
    @{Cell="F26"; Value=18.3919207588041},
    @{Cell="E27"; Value=12.351},
    @{Cell="F27"; Value=10.6078511915737},
    @{Cell="I27"; Value=159.443564150564},
    @{Cell="J27"; Value=0.457178674042294},
    @{Cell="E28"; Value=17.112},
    @{Cell="F28"; Value=14.3089732837144},
    @{Cell="I28"; Value=103.229225339239},
    @{Cell="J28"; Value=0.490301631631449},
    @{Cell="E29"; Value=1.472},
    @{Cell="F29"; Value=1.20701781039144},
    @{Cell="E30"; Value=2.553},
    @{Cell="F30"; Value=2.19837723130137},
    @{Cell="E31"; Value=9.131},
    @{Cell="F31"; Value=6.74770920970798},
    @{Cell="I31"; Value=56.5861720564009},
    @{Cell="J31"; Value=0.839062631602656},
    @{Cell="B33"; Value=191.36},
    @{Cell="E33"; Value=1.15},
    @{Cell="F33"; Value=0.929471728937235},
    @{Cell="I33"; Value=88.3628207784537},
    @{Cell="J33"; Value=0.127192706685308},
    @{Cell="E34"; Value=8.579},
    @{Cell="F34"; Value=6.40726443005721},
    @{Cell="I34"; Value=47.4978920316997},
    @{Cell="J34"; Value=1.66567162730912},
    @{Cell="E35"; Value=13.6965},
    @{Cell="F35"; Value=11.555515223702},
    @{Cell="E36"; Value=4.9565},
    @{Cell="F36"; Value=4.00536243063978}
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

Write-Output "Applied $($updates.Count) cell updates"
